$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Bad Drivers" table -------------------------------------------------
# Remove the retired "Intel(R) Wi-Fi 6 AX201 160MHz - 22.250.1.2" row.
# This shifts the remaining rows (Realtek 155.1, Realtek 8821CE, Totals)
# up by one (old rows 4-6 become new rows 3-5).
$ws.Rows(3).Delete()

# Updated weekly figures for the two remaining drivers.
$ws.Range("C3").Value = 35020
$ws.Range("C4").Value = 398
$ws.Range("D4").Value = 98.3

# Updated totals row (now row 5).
$ws.Range("B5").Value = 26
$ws.Range("C5").Value = 35418

# --- "Good Drivers" table -------------------------------------------------
# Remove the retired "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4" row
# (now at row 13 after the first deletion above).
$ws.Rows(13).Delete()

# Updated sample count for the 6001.15.152.0 driver (now row 14).
$ws.Range("B14").Value = 1033024

# Remove the five trailing "Intel(R) Wi-Fi 6 AX201" rows that are no
# longer part of the Good Drivers list (now rows 21-25).
$ws.Range("A21:A25").EntireRow.Delete()
